$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 129, pushing the existing rows 129-150 down to 130-151.
$ws.Rows("129:129").Insert()

# Populate the newly inserted row 129 with the new data record.
$ws.Range("A129").Value = 5
$ws.Range("B129").Value = "Macroferia Regional de Talca"
$ws.Range("C129").Value = "Maule"
$ws.Range("D129").Value = 44995
$ws.Range("E129").Value = 7
$ws.Range("F129").Value = 100112001
$ws.Range("G129").Value = "Berenjena"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 200
$ws.Range("K129").Value = 8000
$ws.Range("L129").Value = 8000
$ws.Range("M129").Value = 8000
$ws.Range("N129").Value = "$/caja 50 unidades"
$ws.Range("O129").Value = "Región del Maule"
$ws.Range("P129").Value = 160
$ws.Range("Q129").Value = 50
$ws.Range("R129").Value = "Hortaliza"
